$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab from "Neurology" to "Session"
$ws.Name = "Session"

# Delete the last 3 data rows (76, 77, 78) so the used range shrinks to A1:F75
$ws.Rows("76:78").Delete()
